{"js": "// Apply each old->new text replacement for the math worksheet.\n// Each original value is unique in the document, so an exact-match\n// search-and-replace is safe and order-independent.\nconst replacements = [\n  [\"2025-05-31 Saturday\", \"2025-06-01 Sunday\"],\n  [\"23+68=91\", \"98-69=29\"],\n  [\"11+29=40\", \"28+70=98\"],\n  [\"33-3=30\", \"95-84=11\"],\n  [\"25+28=53\", \"96-57=39\"],\n  [\"57+9=66\", \"70-43=27\"],\n  [\"11+31=42\", \"91-71=20\"],\n  [\"60-39=21\", \"36+34=70\"],\n  [\"56+40=96\", \"61-49=12\"],\n  [\"24+7=31\", \"3+10=13\"],\n  [\"39+29=68\", \"95-85=10\"],\n  [\"53-10=43\", \"36+15=51\"],\n  [\"73-36=37\", \"57-43=14\"],\n  [\"70-68=2\", \"28+8=36\"],\n  [\"93-70=23\", \"81-0=81\"],\n  [\"19-1=18\", \"63+29=92\"],\n  [\"54-5=49\", \"84-62=22\"],\n  [\"60+14=74\", \"77-0=77\"],\n  [\"65-54=11\", \"9+82=91\"],\n  [\"39-13=26\", \"30+32=62\"],\n  [\"87+3=90\", \"18+51=69\"],\n  [\"91-44=47\", \"13+57=70\"],\n  [\"52+25=77\", \"22+50=72\"],\n  [\"37+16=53\", \"70-35=35\"],\n  [\"57-28=29\", \"47+13=60\"],\n  [\"8+24=32\", \"56+38=94\"],\n  [\"93-20=73\", \"24+36=60\"],\n  [\"42-9=33\", \"34+63=97\"],\n  [\"56+39=95\", \"3+87=90\"],\n  [\"96-31=65\", \"51-23=28\"],\n  [\"40-37=3\", \"64-25=39\"],\n  [\"78-15=63\", \"24+51=75\"],\n  [\"18-12=6\", \"16+18=34\"],\n  [\"36+3=39\", \"54-30=24\"],\n  [\"5+13=18\", \"53+0=53\"],\n  [\"56-46=10\", \"34-20=14\"],\n  [\"38-35=3\", \"42+7=49\"],\n  [\"68-41=27\", \"54+5=59\"],\n  [\"71-44=27\", \"97-70=27\"],\n  [\"36+58=94\", \"68-29=39\"],\n  [\"12+0=12\", \"9+17=26\"],\n  [\"11+0=11\", \"13+23=36\"],\n  [\"15+56=71\", \"59+19=78\"],\n  [\"24-12=12\", \"60-40=20\"],\n  [\"1+79=80\", \"63+33=96\"],\n  [\"50-1=49\", \"80+1=81\"],\n  [\"58-33=25\", \"40-33=7\"],\n  [\"66-38=28\", \"64-52=12\"],\n  [\"48+13=61\", \"85-71=14\"],\n  [\"38+2=40\", \"79-41=38\"],\n  [\"38+27=65\", \"9+60=69\"],\n  [\"41+34=75\", \"46+14=60\"],\n  [\"65+8=73\", \"34-12=22\"],\n  [\"35-19=16\", \"18+40=58\"],\n  [\"31+62=93\", \"18+28=46\"],\n  [\"76-72=4\", \"22+50=72\"],\n  [\"9+27=36\", \"85-74=11\"],\n  [\"84-40=44\", \"25+73=98\"],\n  [\"42+47=89\", \"16+17=33\"],\n  [\"70-31=39\", \"74-53=21\"],\n  [\"53-34=19\", \"61-50=11\"],\n  [\"36+55=91\", \"11+17=28\"],\n  [\"64-50=14\", \"55-17=38\"],\n  [\"91-68=23\", \"9-6=3\"],\n  [\"9+40=49\", \"35+36=71\"],\n  [\"89-46=43\", \"70+28=98\"],\n  [\"71-62=9\", \"91-38=53\"],\n  [\"20+36=56\", \"22-0=22\"],\n  [\"76+9=85\", \"76-20=56\"],\n  [\"15+66=81\", \"28-27=1\"],\n  [\"57+35=92\", \"15-6=9\"],\n  [\"11+18=29\", \"36-9=27\"],\n  [\"51-26=25\", \"40+17=57\"],\n  [\"19+21=40\", \"14+18=32\"],\n  [\"5+87=92\", \"83-81=2\"],\n  [\"64-18=46\", \"6+24=30\"],\n  [\"21+60=81\", \"88+3=91\"],\n  [\"87+10=97\", \"60+5=65\"],\n  [\"17+11=28\", \"42+18=60\"],\n  [\"79-66=13\", \"16-8=8\"],\n  [\"17+67=84\", \"56-29=27\"],\n  [\"97-84=13\", \"69-6=63\"],\n  [\"91-49=42\", \"39-37=2\"],\n  [\"25+1=26\", \"76-46=30\"],\n  [\"86-8=78\", \"71+8=79\"],\n  [\"45+43=88\", \"72+25=97\"],\n  [\"85-7=78\", \"31+68=99\"],\n  [\"4+93=97\", \"84-69=15\"],\n  [\"11+70=81\", \"97-38=59\"],\n  [\"85-45=40\", \"97-26=71\"],\n  [\"69-46=23\", \"20+59=79\"],\n  [\"84-2=82\", \"86+8=94\"],\n  [\"18-8=10\", \"71-50=21\"],\n  [\"58-34=24\", \"38+43=81\"],\n  [\"49+20=69\", \"84-34=50\"],\n  [\"66-10=56\", \"58-10=48\"],\n  [\"72-18=54\", \"33+9=42\"],\n  [\"39+60=99\", \"53+11=64\"],\n  [\"67-65=2\", \"20+0=20\"],\n  [\"12+67=79\", \"18+80=98\"],\n  [\"30+41=71\", \"53-1=52\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply each old->new text replacement for the math worksheet.\n# Each original value is unique in the document, so an exact-match\n# Find/Replace (MatchWholeWord off, exact text) is safe and order-independent.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2025-05-31 Saturday\", \"2025-06-01 Sunday\"),\n  @(\"23+68=91\", \"98-69=29\"),\n  @(\"11+29=40\", \"28+70=98\"),\n  @(\"33-3=30\", \"95-84=11\"),\n  @(\"25+28=53\", \"96-57=39\"),\n  @(\"57+9=66\", \"70-43=27\"),\n  @(\"11+31=42\", \"91-71=20\"),\n  @(\"60-39=21\", \"36+34=70\"),\n  @(\"56+40=96\", \"61-49=12\"),\n  @(\"24+7=31\", \"3+10=13\"),\n  @(\"39+29=68\", \"95-85=10\"),\n  @(\"53-10=43\", \"36+15=51\"),\n  @(\"73-36=37\", \"57-43=14\"),\n  @(\"70-68=2\", \"28+8=36\"),\n  @(\"93-70=23\", \"81-0=81\"),\n  @(\"19-1=18\", \"63+29=92\"),\n  @(\"54-5=49\", \"84-62=22\"),\n  @(\"60+14=74\", \"77-0=77\"),\n  @(\"65-54=11\", \"9+82=91\"),\n  @(\"39-13=26\", \"30+32=62\"),\n  @(\"87+3=90\", \"18+51=69\"),\n  @(\"91-44=47\", \"13+57=70\"),\n  @(\"52+25=77\", \"22+50=72\"),\n  @(\"37+16=53\", \"70-35=35\"),\n  @(\"57-28=29\", \"47+13=60\"),\n  @(\"8+24=32\", \"56+38=94\"),\n  @(\"93-20=73\", \"24+36=60\"),\n  @(\"42-9=33\", \"34+63=97\"),\n  @(\"56+39=95\", \"3+87=90\"),\n  @(\"96-31=65\", \"51-23=28\"),\n  @(\"40-37=3\", \"64-25=39\"),\n  @(\"78-15=63\", \"24+51=75\"),\n  @(\"18-12=6\", \"16+18=34\"),\n  @(\"36+3=39\", \"54-30=24\"),\n  @(\"5+13=18\", \"53+0=53\"),\n  @(\"56-46=10\", \"34-20=14\"),\n  @(\"38-35=3\", \"42+7=49\"),\n  @(\"68-41=27\", \"54+5=59\"),\n  @(\"71-44=27\", \"97-70=27\"),\n  @(\"36+58=94\", \"68-29=39\"),\n  @(\"12+0=12\", \"9+17=26\"),\n  @(\"11+0=11\", \"13+23=36\"),\n  @(\"15+56=71\", \"59+19=78\"),\n  @(\"24-12=12\", \"60-40=20\"),\n  @(\"1+79=80\", \"63+33=96\"),\n  @(\"50-1=49\", \"80+1=81\"),\n  @(\"58-33=25\", \"40-33=7\"),\n  @(\"66-38=28\", \"64-52=12\"),\n  @(\"48+13=61\", \"85-71=14\"),\n  @(\"38+2=40\", \"79-41=38\"),\n  @(\"38+27=65\", \"9+60=69\"),\n  @(\"41+34=75\", \"46+14=60\"),\n  @(\"65+8=73\", \"34-12=22\"),\n  @(\"35-19=16\", \"18+40=58\"),\n  @(\"31+62=93\", \"18+28=46\"),\n  @(\"76-72=4\", \"22+50=72\"),\n  @(\"9+27=36\", \"85-74=11\"),\n  @(\"84-40=44\", \"25+73=98\"),\n  @(\"42+47=89\", \"16+17=33\"),\n  @(\"70-31=39\", \"74-53=21\"),\n  @(\"53-34=19\", \"61-50=11\"),\n  @(\"36+55=91\", \"11+17=28\"),\n  @(\"64-50=14\", \"55-17=38\"),\n  @(\"91-68=23\", \"9-6=3\"),\n  @(\"9+40=49\", \"35+36=71\"),\n  @(\"89-46=43\", \"70+28=98\"),\n  @(\"71-62=9\", \"91-38=53\"),\n  @(\"20+36=56\", \"22-0=22\"),\n  @(\"76+9=85\", \"76-20=56\"),\n  @(\"15+66=81\", \"28-27=1\"),\n  @(\"57+35=92\", \"15-6=9\"),\n  @(\"11+18=29\", \"36-9=27\"),\n  @(\"51-26=25\", \"40+17=57\"),\n  @(\"19+21=40\", \"14+18=32\"),\n  @(\"5+87=92\", \"83-81=2\"),\n  @(\"64-18=46\", \"6+24=30\"),\n  @(\"21+60=81\", \"88+3=91\"),\n  @(\"87+10=97\", \"60+5=65\"),\n  @(\"17+11=28\", \"42+18=60\"),\n  @(\"79-66=13\", \"16-8=8\"),\n  @(\"17+67=84\", \"56-29=27\"),\n  @(\"97-84=13\", \"69-6=63\"),\n  @(\"91-49=42\", \"39-37=2\"),\n  @(\"25+1=26\", \"76-46=30\"),\n  @(\"86-8=78\", \"71+8=79\"),\n  @(\"45+43=88\", \"72+25=97\"),\n  @(\"85-7=78\", \"31+68=99\"),\n  @(\"4+93=97\", \"84-69=15\"),\n  @(\"11+70=81\", \"97-38=59\"),\n  @(\"85-45=40\", \"97-26=71\"),\n  @(\"69-46=23\", \"20+59=79\"),\n  @(\"84-2=82\", \"86+8=94\"),\n  @(\"18-8=10\", \"71-50=21\"),\n  @(\"58-34=24\", \"38+43=81\"),\n  @(\"49+20=69\", \"84-34=50\"),\n  @(\"66-10=56\", \"58-10=48\"),\n  @(\"72-18=54\", \"33+9=42\"),\n  @(\"39+60=99\", \"53+11=64\"),\n  @(\"67-65=2\", \"20+0=20\"),\n  @(\"12+67=79\", \"18+80=98\"),\n  @(\"30+41=71\", \"53-1=52\"),\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
